$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.109.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.653.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2606'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07801'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.669.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.880.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5477'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.099.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.586'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.013'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.17%  '
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.212'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.459'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05773'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.24%  '
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.546'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.263'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.595'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.416'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9463'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5755'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01612'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8536'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.715'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.031.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.794.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.90'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.002'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4333'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05143'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.826'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('E51').Value = '  -1.51%  '
